$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Text (shared-string) edits -----------------------------------------

# F2 / "HIPPo" description: "..." -> "…", drop trailing space before the
# line break, and wrap "Variantes :" in <u>...</u>
$ws.Range("F2").Value = @'
HIPPo : C’est le chef (la personne la plus importante, la mieux payée, la plus forte…) qui décide pour tout le monde.
<u>Variantes :</u>
Oligarchie (les plus âgés décident). 
Vote féminin/masculin (seul.e.s les filles/garçons décident). 
'@

# F6 / "Consentement": drop the "une " before "majorité"
$ws.Range("F6").Value = @'
Les participants peuvent répondre par 3 choix :
 - Je suis d’accord : <b>pouce levé</b>
 - Je consens : <b>main horizontale ou à plat</b>. En consentant, j’exprime une opinion. Je consens que la personne ayant proposé aille plus loin.
 - Je mets mon veto : <b>pouce en bas</b>. Si je mets un veto, je dois expliquer les raisons de mon veto et proposer une porte de sortie. Je serais d’accord si… 
Pour que la décision soit prise, il ne faut aucun veto et majorité de ‘Pour’.
'@

# F9 / "Choix commun": reword a few phrases
$ws.Range("F9").Value = @'
Chaque votant donne à chacune des propositions une mention claire pour tous : Excellent/Très bien/Bien/Assez bien/Passable/Insuffisant/À rejeter. 
En prenant en compte l’ensemble des votes, cela donne à chaque proposition un profil. 
La proposition, qui a le meilleur profil est choisie : c’est celui qui aura les mentions les plus favorables. 
Si aucune proposition n’a de profil très favorable, on peut rediscuter celles-ci !
'@

# F12 / "Nappes tournantes": reword + underline "Astuces :" + re-wrap last line
$ws.Range("F12").Value = @'
Une ou plusieurs problématiques sont inscrites sur <b>chaque coin d’une nappe</b>.
Les participants doivent répondre, en 3-5 min, à la question <b>devant eux</b>.
A la fin du temps imparti, <b>la nappe tourne</b> devant les participants. Et on recommence.
<u>Astuces :</u>
Faire autant de tours que de participants. Chacun pourra ainsi enrichir chaque sujet.
Augmenter le temps imparti si les participants n’ont pas le temps de 
tout lire.
'@

# F14 / "Mind map": "ilôt" -> "ilôts"
$ws.Range("F14").Value = @'
Au centre de la feuille/tableau, Écrire <b>un</b> mot/une phrase qui résume la problématique.
Raccrocher les mots/concepts qui <b>gravitent</b> autour, les mettre sous forme d'ilôts.
Enrichir les concepts représentés.
'@

# --- 2. Row heights (re-flowed after the text/format pass) -----------------

$ws.Rows.Item(2).RowHeight = 79.45
$ws.Rows.Item(6).RowHeight = 98.95
$ws.Rows.Item(8).RowHeight = 33.7
$ws.Rows.Item(9).RowHeight = 89.2
$ws.Rows.Item(12).RowHeight = 118.45
$ws.Rows.Item(14).RowHeight = 50.2

# --- 3. Row 15 clean-up: drop the two truly-empty helper cells -------------

$ws.Range("C15").ClearContents()
$ws.Range("E15").ClearContents()

# --- 4. New "poker format" helper cell in G6 (Calibri font, no content) ----

$ws.Range("G6").Font.Name = "Calibri"
$ws.Range("G6").Font.Size = 10

# --- 5. Header / footer font: "...,Regular" -> "...,Book" ------------------

$ws.PageSetup.CenterHeader = '&"DejaVu Serif,Book"&12&A'
$ws.PageSetup.CenterFooter = '&"DejaVu Serif,Book"&12Page &P'

# --- 6. View state: active cell / scroll position ---------------------------

$ws.Range("F13").Select()
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
